$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (repayment_20250901_20250918 (1) -> (2))
$ws.Name = "repayment_20250901_20250918 (2)"

function Set-TextValue($range, $text) {
    # Force a genuinely text-typed cell (matches shared-string cells in the
    # source file) instead of letting Excel auto-convert numeric-looking
    # strings (e.g. "31,414,185.00") into a real number. ClearFormats()
    # afterwards drops the temporary "@" number format again so the cell
    # keeps using the default style, just like the original file.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2 - Debora Retima Sihombing
$ws.Range("H2").Value = 19.483000000000001

# Row 3 - Romli
$ws.Range("H3").Value = 25.158999999999999

# Row 4 - Aldi Taufik
$ws.Range("H4").Value = 15.428000000000001

# Row 5 - Yandi Nugraha
$ws.Range("D5").Value = 51
Set-TextValue $ws.Range("E5") "31,414,185.00"
Set-TextValue $ws.Range("G5") "9.33"
$ws.Range("H5").Value = 22.131

# Row 6 - Axl Wicaksono
$ws.Range("D6").Value = 44
Set-TextValue $ws.Range("E6") "29,021,986.00"
Set-TextValue $ws.Range("G6") "8.28"
$ws.Range("H6").Value = 16.190999999999999

# Row 7 - Riska Nurlita
$ws.Range("H7").Value = 12.834

# Row 8 - Annisa Putri Restu
Set-TextValue $ws.Range("E8") "41,197,087.00"
Set-TextValue $ws.Range("G8") "11.42"
$ws.Range("H8").Value = 25.597999999999999

# Row 9 - Azizah Rahmawati
$ws.Range("H9").Value = 13.742000000000001

# Row 10 - Erlangga Hutama
$ws.Range("H10").Value = 14.651

# Row 11 - Erick Ervan Dewanggga
$ws.Range("D11").Value = 51
Set-TextValue $ws.Range("E11") "44,106,636.00"
Set-TextValue $ws.Range("G11") "12.13"
$ws.Range("H11").Value = 14.093

# Row 12 - Ridhoi Berkat Zebua
$ws.Range("H12").Value = 20.692

# Row 13 - Fadilah Damayanti
$ws.Range("H13").Value = 17.847000000000001

# Row 15 - Adistira Winditya P
$ws.Range("H15").Value = 12.191000000000001

# Row 16 - Sucika Wardani
$ws.Range("D16").Value = 47
Set-TextValue $ws.Range("E16") "32,176,238.00"
Set-TextValue $ws.Range("G16") "8.87"
$ws.Range("H16").Value = 11.395

# Row 17 - Wasti Feronika Sihombing
$ws.Range("D17").Value = 46
Set-TextValue $ws.Range("E17") "33,338,944.00"
Set-TextValue $ws.Range("G17") "9.54"
$ws.Range("H17").Value = 21.704000000000001

# Row 18 - Nuraini
$ws.Range("H18").Value = 12.416
